# Refresh the crypto price/volume table (cryptos.xlsx, sheet1) with the
# latest scraped values. Column D ("Price") and E ("Volume(1h)") are stored
# as plain text in the source data (coinranking.com scrape), and row 48/49
# swap which coin (NEARProtocol / PaxDollar) occupies that rank slot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.888.73'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '1.815.75'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'309.90"
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = "'0.4684"
$ws.Range('E7').Value = '  +1.31%  '
$ws.Range('D8').Value = "'0.3695"
$ws.Range('E8').Value = '  -1.64%  '
$ws.Range('D9').Value = "'0.07377"
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').Value = "'0.8706"
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').Value = '1.839.68'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').Value = "'5.367"
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = "'92.24"
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = "'0.07076"
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = "'6.502"
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = "'0.000008723"
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('D21').Value = '26.906.56'
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('D22').Value = "'5.334"
$ws.Range('E22').Value = '  +0.42%  '
$ws.Range('D23').Value = "'10.54"
$ws.Range('E23').Value = '  -2.99%  '
$ws.Range('D24').Value = '2.020.43'
$ws.Range('E24').Value = '  -1.12%  '
$ws.Range('D25').Value = "'1.892"
$ws.Range('E25').Value = '  -1.70%  '
$ws.Range('D26').Value = "'151.88"
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').Value = "'2.194"
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').Value = "'18.37"
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').Value = "'5.303"
$ws.Range('E29').Value = '  +0.81%  '
$ws.Range('D30').Value = "'115.47"
$ws.Range('E30').Value = '  -1.24%  '
$ws.Range('D31').Value = "'0.08922"
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('D32').Value = "'0.7667"
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('D33').Value = "'1.161"
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('D34').Value = "'4.481"
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('D35').Value = "'2.921"
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = "'1.097"
$ws.Range('E37').Value = '  -2.77%  '
$ws.Range('D38').Value = "'0.01960"
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').Value = "'2.941"
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').Value = "'7.245"
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').Value = "'0.5340"
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').Value = "'2.345"
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('D44').Value = "'0.1662"
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('D45').Value = "'8.448"
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('D46').Value = "'0.4941"
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('D47').Value = "'10.44"
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Value = "'1.000"
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'1.671"
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').Value = "'102.79"
$ws.Range('E50').Value = '  -2.05%  '
$ws.Range('D51').Value = "'0.06283"
$ws.Range('E51').Value = '  -0.61%  '
